$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price cells so Excel
# does not auto-convert them to actual numbers (loses formatting).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values.
$ws.Range("D2").Value = '69.687.71'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '3.515.25'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '616.62'
$ws.Range("E5").Value = '  +5.56%  '
$ws.Range("D6").Value = '192.08'
$ws.Range("E6").Value = '  +1.65%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.212'
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("D10").Value = '0.655'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '53.47'
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("D12").Value = '0.0000309'
$ws.Range("E12").Value = '  -2.67%  '
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '4.077.29'
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").Value = '616.82'
$ws.Range("E15").Value = '  +7.58%  '
$ws.Range("D16").Value = '69.732.44'
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("D17").Value = '19.01'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '12.59'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '3.516.70'
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '0.989'
$ws.Range("E21").Value = '  -1.57%  '
$ws.Range("D22").Value = '17.17'
$ws.Range("E22").Value = '  -2.48%  '
$ws.Range("D23").Value = '106.94'
$ws.Range("E23").Value = '  +13.50%  '
$ws.Range("D24").Value = '4.72'
$ws.Range("E24").Value = '  +3.07%  '
$ws.Range("D25").Value = '5.08'
$ws.Range("E25").Value = '  +3.43%  '
$ws.Range("D26").Value = '3.08'
$ws.Range("E26").Value = '  +5.03%  '
$ws.Range("D27").Value = '11.01'
$ws.Range("E27").Value = '  -1.75%  '
$ws.Range("D28").Value = '9.73'
$ws.Range("E28").Value = '  +4.68%  '
$ws.Range("D29").Value = '34.02'
$ws.Range("E29").Value = '  +4.05%  '
$ws.Range("E30").Value = '  -3.07%  '
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("D32").Value = '3.94'
$ws.Range("E32").Value = '  +3.39%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = '63.47'
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("D35").Value = '3.13'
$ws.Range("E35").Value = '  -5.62%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = '3.660.66'
$ws.Range("E37").Value = '  +1.12%  '
$ws.Range("D38").Value = '514.80'
$ws.Range("E38").Value = '  -3.04%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '3.63'
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '0.395'
$ws.Range("E40").Value = '  -3.99%  '
$ws.Range("D41").Value = '0.0₃0790'
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("D42").Value = '36.74'
$ws.Range("E42").Value = '  -4.43%  '
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").Value = '0.0465'
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("E46").Value = '  +3.41%  '
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  -3.92%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '8.76'
$ws.Range("E48").Value = '  -5.72%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("D50").Value = '132.20'
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("E51").Value = '  -6.79%  '
